$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "基金受益憑證" (fund) sheet: turn the duplicate-data row 1 into a proper
#    header row of field names, and append the standard trailing metadata
#    columns (property_category .. index) that every other sheet already has.
# ---------------------------------------------------------------------------
$wsFund = $wb.Worksheets.Item("基金受益憑證")

$wsFund.Range("B1").Value = "name"
$wsFund.Range("C1").Value = "owner"
$wsFund.Range("D1").Value = "dealer"
$wsFund.Range("E1").Value = "quantity"
$wsFund.Range("F1").Value = "face_value"
$wsFund.Range("G1").Value = "currency"
$wsFund.Range("H1").Value = "total"
$wsFund.Range("I1").Value = "property_category"
$wsFund.Range("J1").Value = "category"
$wsFund.Range("K1").Value = "date"
$wsFund.Range("L1").Value = "legislator_name"
$wsFund.Range("M1").Value = "legislator_id"
$wsFund.Range("N1").Value = "source_file"
$wsFund.Range("O1").Value = "index"

# Dates need to stay plain text ("2011-12-26"), not get reinterpreted as a
# date serial by Excel's auto-detection, so force text format first.
$wsFund.Range("K2:K3").NumberFormat = "@"

$wsFund.Range("I2").Value = "fund"
$wsFund.Range("J2").Value = "normal"
$wsFund.Range("K2").Value = "2011-12-26"
$wsFund.Range("L2").Value = "黃昭順"
$wsFund.Range("M2").Value = 665
$wsFund.Range("N2").Value = "tmp43441"
$wsFund.Range("O2").Value = 84

$wsFund.Range("I3").Value = "fund"
$wsFund.Range("J3").Value = "normal"
$wsFund.Range("K3").Value = "2011-12-26"
$wsFund.Range("L3").Value = "黃昭順"
$wsFund.Range("M3").Value = 665
$wsFund.Range("N3").Value = "tmp43441"
$wsFund.Range("O3").Value = 85

# ---------------------------------------------------------------------------
# 2) Drop the "其他有價證券" sheet entirely (its data was junk/mis-parsed
#    header fragments, not real records) so "保險" shifts up to take its
#    place in the tab order.
# ---------------------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("其他有價證券")
$wsOther.Delete()
